# Update "Förändrad" date column (C) from 2025-06-19 to 2025-06-20
# for all data rows (rows 2 through 43) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45827) {
        $cell.Value = 45828
    }
}
